# Fixed some bugs in CalcAdjacentPay
#
# The unit-test data table (columns: symbol, reel1, reel2, reel3, reel4,
# reel5) had a batch of rows shuffled into the wrong order. This restores
# the correct per-row contents for the affected rows; row 1 (headers),
# the rows that were already correct, and the totals row (26) are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowNumbers = @(2, 3, 4, 5, 6, 7, 9, 10, 11, 12, 13, 14, 15, 16, 19, 20, 21)

$rowData = @(
    @(902, 1, 0, 0, 0, 0),
    @(501, 9, 52, 30, 75, 45),
    @(201, 9, 30, 15, 45, 30),
    @(801, 3, 67, 65, 52, 45),
    @(1203, 3, 15, 15, 15, 15),
    @(901, 16, 15, 45, 60, 60),
    @(701, 3, 90, 45, 97, 15),
    @(101, 9, 30, 15, 60, 15),
    @(401, 9, 48, 67, 75, 45),
    @(1202, 2, 10, 10, 10, 10),
    @(1001, 18, 30, 75, 60, 72),
    @(601, 9, 60, 67, 60, 42),
    @(1201, 2, 10, 10, 10, 10),
    @(1, 0, 2, 2, 2, 2),
    @(3, 0, 3, 3, 3, 3),
    @(802, 0, 4, 5, 4, 0),
    @(1101, 0, 15, 30, 30, 0)
)

for ($idx = 0; $idx -lt $rowNumbers.Length; $idx++) {
    $r = $rowNumbers[$idx]
    $vals = $rowData[$idx]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
}

Write-Output "CalcAdjacentPay test data rows corrected."
